# Applies the "big stimulus update" to the stimulus list:
#  - rename the "face" image category to "book" (face//face_NN.jpg -> book//book_NN.jpg)
#  - spell out the abbreviated correct_ans codes in column L (r -> right, y -> left, b -> center)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

$ansMap = @{ "r" = "right"; "y" = "left"; "b" = "center" }

for ($row = 2; $row -le $lastRow; $row++) {

    # columns A-D hold the stimulus filenames (promptFile, correctFile, dist_01File, dist_02File)
    for ($col = 1; $col -le 4; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $cur = $cell.Value2
        if ($cur -like "face//face_*") {
            $cell.Value2 = $cur -replace "face//face_", "book//book_"
        }
    }

    # column L holds the abbreviated correct answer code
    $ansCell = $ws.Cells.Item($row, 12)
    $curAns = $ansCell.Value2
    if ($ansMap.ContainsKey($curAns)) {
        $ansCell.Value2 = $ansMap[$curAns]
    }
}
